$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks from the product URL cells (C2, C3) - these
# pointed at www.prod1.com / www.prod2.com and are being replaced by local
# image paths used from wwwroot.
$null = $ws.Hyperlinks.Delete()

# Replace the URL text with paths to the working images copied into wwwroot.
$ws.Range("C2").Value = "./images/dota.jpg"
$ws.Range("C3").Value = "./images/dotalol.jpg"

# Update the active selection to match the saved view from the author's session.
$null = $ws.Range("I9").Select()
